$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking string to be stored as text, matching the
# original inline-string cell type, without leaving a lingering cell style.
function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.ClearFormats()
}

$ws.Range('D2').Value = '36.085.89'
$ws.Range('E2').Value = '  -3.75%  '
$ws.Range('D3').Value = '1.963.23'
$ws.Range('E3').Value = '  -4.08%  '
$ws.Range('E4').Value = '  +0.16%  '
Set-TextValue 'D5' '241.86'
$ws.Range('E5').Value = '  -4.17%  '
$ws.Range('E6').Value = '  -3.23%  '
Set-TextValue 'D7' '62.06'
$ws.Range('E7').Value = '  -6.26%  '
$ws.Range('E8').Value = '  +0.07%  '
Set-TextValue 'D9' '0.373'
$ws.Range('E9').Value = '  -1.05%  '
Set-TextValue 'D10' '56.15'
$ws.Range('E10').Value = '  -5.62%  '
Set-TextValue 'D11' '0.0808'
$ws.Range('E11').Value = '  +6.98%  '
$ws.Range('E12').Value = '  -1.06%  '
Set-TextValue 'D13' '0.855'
$ws.Range('E13').Value = '  -6.12%  '
Set-TextValue 'D14' '21.99'
$ws.Range('E14').Value = '  +6.04%  '
Set-TextValue 'D15' '14.01'
$ws.Range('E15').Value = '  -7.87%  '
$ws.Range('D16').Value = '2.250.17'
$ws.Range('E16').Value = '  -4.07%  '
Set-TextValue 'D17' '5.41'
$ws.Range('E17').Value = '  -3.90%  '
$ws.Range('D18').Value = '1.966.38'
$ws.Range('E18').Value = '  -4.01%  '
$ws.Range('D19').Value = '35.982.95'
$ws.Range('E19').Value = '  -3.79%  '
Set-TextValue 'D20' '71.05'
$ws.Range('E20').Value = '  -3.35%  '
$ws.Range('D21').Value = '0.0₃0854'
$ws.Range('E21').Value = '  -2.96%  '
Set-TextValue 'D22' '236.77'
$ws.Range('E22').Value = '  -0.34%  '
$ws.Range('E23').Value = '  -3.08%  '
$ws.Range('E24').Value = '  +0.19%  '
$ws.Range('E25').Value = '  -7.14%  '
$ws.Range('E26').Value = '  -3.37%  '
Set-TextValue 'D27' '9.74'
$ws.Range('E27').Value = '  +1.45%  '
Set-TextValue 'D28' '159.95'
$ws.Range('E28').Value = '  -3.39%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D29' '0.134'
$ws.Range('E29').Value = '  +22.09%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D30' '19.83'
$ws.Range('E30').Value = '  -0.38%  '
$ws.Range('E31').Value = '  -1.98%  '
$ws.Range('E32').Value = '  -6.60%  '
$ws.Range('E33').Value = '  -8.06%  '
Set-TextValue 'D34' '0.0619'
$ws.Range('E34').Value = '  +0.77%  '
Set-TextValue 'D35' '4.41'
$ws.Range('E35').Value = '  -6.90%  '
Set-TextValue 'D36' '6.30'
$ws.Range('E36').Value = '  +4.22%  '
$ws.Range('E37').Value = '  +0.18%  '
Set-TextValue 'D38' '2.27'
$ws.Range('E38').Value = '  -7.55%  '
Set-TextValue 'D39' '1.83'
$ws.Range('E39').Value = '  +0.92%  '
Set-TextValue 'D40' '3.06'
$ws.Range('E40').Value = '  +12.63%  '
Set-TextValue 'D41' '0.0984'
$ws.Range('E41').Value = '  -4.99%  '
$ws.Range('E42').Value = '  -1.79%  '
$ws.Range('E43').Value = '  -3.51%  '
$ws.Range('E44').Value = '  -4.58%  '
Set-TextValue 'D45' '1.08'
$ws.Range('E45').Value = '  -5.10%  '
Set-TextValue 'D46' '92.18'
$ws.Range('E46').Value = '  -3.52%  '
Set-TextValue 'D47' '15.99'
$ws.Range('E47').Value = '  -6.24%  '
Set-TextValue 'D48' '7.54'
$ws.Range('E48').Value = '  -6.95%  '
$ws.Range('D49').Value = '1.336.55'
$ws.Range('E49').Value = '  -6.53%  '
$ws.Range('E50').Value = '  -5.03%  '
$ws.Range('D51').Value = '2.143.03'
$ws.Range('E51').Value = '  -4.14%  '
